# =====================================================================
# Update "Sample Section" template workbook (sample-section.xlsx)
#  1. Refresh all column-header cell comments (A1:V1) with new wording
#  2. Add "Water" to the storage_medium lookup list and re-sort entries
#  3. Point the storage_medium data-validation at the now-larger range
#  4. Bump the pav:createdOn timestamp on the .metadata sheet
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Update cell comments on "Sample Section" sheet (row-1 headers) ----
$ws = $wb.Worksheets.Item("Sample Section")

$c_A1 = @"
(Required) The unique identifier from HuBMAP or SenNet for the source (parent
data) from which the sample was derived. Example: HBM122.EFGH.789
"@
$ws.Range("A1").Comment.Text($c_A1)

$c_B1 = @"
(Required) The unique HuBMAP or SenNet identifier assigned to the sample by the
ingest portal. Example: HBM743.CKJW.876
"@
$ws.Range("B1").Comment.Text($c_B1)

$c_C1 = @"
A locally assigned identifier provided by the data provider for the dataset. It
is used to reference an external metadata record that may be maintained
independently, enabling traceability and supporting provenance tracking.
Example: Visium_9OLC_A4_S1
"@
$ws.Range("C1").Comment.Text($c_C1)

$c_D1 = @"
(Required) The DOI for the protocols.io page that details the assay or the
procedures used for sample procurement and preparation. For example, in the case
of an imaging assay, the protocol may start with tissue section staining and end
with the generation of an OME-TIFF file. The documented protocol should also
include any image processing steps involved in producing the final OME-TIFF.
Example: https://dx.doi.org/10.17504/protocols.io.eq2lyno9qvx9/v1
"@
$ws.Range("D1").Comment.Text($c_D1)

$c_E1 = @"
(Required) The length of time the sample was stored prior to processing it. For
assays performed on tissue sections, this refers to how long the tissue section
(e.g., slide) was stored before the assay began (e.g., imaging). For assays
performed on suspensions, such as sequencing, it refers to how long the
suspension was stored before library construction started. Example: 12
"@
$ws.Range("E1").Comment.Text($c_E1)

$c_F1 = @"
(Required) The unit of measurement used to specify the source storage duration
value. Example: hour
"@
$ws.Range("F1").Comment.Text($c_F1)

$c_G1 = @"
(Required) The medium used during the sample preparation process. If no specific
medium was utilized, enter "None". If medium was not recorded, enter "Unknown".
Example: Fresh frozen CMC
"@
$ws.Range("G1").Comment.Text($c_G1)

$c_H1 = @"
(Required) The condition under which the sample preparation took place, such as
whether the sample was placed on dry ice during the process. If preparation
condition was not recorded, enter "Unknown". Example: Frozen on dry ice
"@
$ws.Range("H1").Comment.Text($c_H1)

$c_I1 = @"
The duration for which the tissue was handled prior to its initial preservation.
Example: 120
"@
$ws.Range("I1").Comment.Text($c_I1)

$c_J1 = @"
The unit of measurement for the processing time value. If processing time is not
specified, this field may be left blank. Example: minute
"@
$ws.Range("J1").Comment.Text($c_J1)

$c_K1 = @"
(Required) The medium used to preserve the sample. If no specific medium was
utilized, enter "None". If medium was not recorded, enter "Unknown". Example:
FFPE (Paraffin embedded)
"@
$ws.Range("K1").Comment.Text($c_K1)

$c_L1 = @"
(Required) The method used to store the sample after preparation and prior to
performing the assay. If no specific storage method was utilized, enter "None".
If storage method was not recorded, enter "Unknown". Example: Frozen in dry ice
"@
$ws.Range("L1").Comment.Text($c_L1)

$c_M1 = @"
The quality criteria used to assess the sample, which may include metrics such
as RIN (e.g., RIN: 8.7) or visual inspection parameters for suspensions prior to
cell lysis. These criteria can be captured at a high level with general terms
like "OK" or "not OK" or with more specific descriptors such as "debris" "clump"
or "low clump". Example: RIN: 8.7, low clump, no visible debris
"@
$ws.Range("M1").Comment.Text($c_M1)

$c_N1 = @"
The key variables in the histopathological report that are crucial for assessing
the tissue, including the absence of necrosis, comments on tissue composition,
descriptions of significant pathology, and high-level assessments of
inflammation or fibrosis. Example: No necrosis observed; tissue composed
predominantly of hepatocytes with mild portal inflammation and minimal fibrosis
"@
$ws.Range("N1").Comment.Text($c_N1)

$c_O1 = @"
(Required) The thickness of an object in question. Example: 10
"@
$ws.Range("O1").Comment.Text($c_O1)

$c_P1 = @"
(Required) The unit of measurement for the thickness value. If no thickness
measurement is specified, this field may be left blank. Example: mm
"@
$ws.Range("P1").Comment.Text($c_P1)

$c_Q1 = @"
(Required) The index number assigned to the tissue section, with numbering
beginning at 1 for sections within a block. Example: 1
"@
$ws.Range("Q1").Comment.Text($c_Q1)

$c_R1 = @"
The area of the object being measured. Example: 100
"@
$ws.Range("R1").Comment.Text($c_R1)

$c_S1 = @"
The unit of measurement used to define the area. If no area value is specified,
this field may be left blank. Example: mm^2
"@
$ws.Range("S1").Comment.Text($c_S1)

$c_T1 = @"
Indicates whether the section was rehydrated. Example: No
"@
$ws.Range("T1").Comment.Text($c_T1)

$c_U1 = @"
Miscellaneous details about the sample that are not captured in the existing
metadata fields. Example: Sample was stored at 4°C for 48 hours prior to
processing due to equipment maintenance delay
"@
$ws.Range("U1").Comment.Text($c_U1)

$c_V1 = @"
(Required) The unique string identifier for the metadata specification version,
which is easily interpretable by computers for purposes of data validation and
processing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9
"@
$ws.Range("V1").Comment.Text($c_V1)

# ---- 2. Rewrite "storage_medium" lookup sheet (adds "Water", re-sorted) ----
$wsSM = $wb.Worksheets.Item("storage_medium")
$wsSM.Cells.Item(1, 1).Value = 'Water'
$wsSM.Cells.Item(1, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65147'
$wsSM.Cells.Item(2, 1).Value = 'OCT'
$wsSM.Cells.Item(2, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63523'
$wsSM.Cells.Item(3, 1).Value = 'NBF (Neutral Buffered Formalin)'
$wsSM.Cells.Item(3, 2).Value = 'http://purl.obolibrary.org/obo/OBIB_0000213'
$wsSM.Cells.Item(4, 1).Value = 'Allprotect tissue reagent (ALL)'
$wsSM.Cells.Item(4, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118'
$wsSM.Cells.Item(5, 1).Value = 'DMSO (no serum)'
$wsSM.Cells.Item(5, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000115'
$wsSM.Cells.Item(6, 1).Value = 'PFA (Paraformaldehyde)'
$wsSM.Cells.Item(6, 2).Value = 'http://purl.obolibrary.org/obo/CHEBI_61538'
$wsSM.Cells.Item(7, 1).Value = 'Unknown'
$wsSM.Cells.Item(7, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998'
$wsSM.Cells.Item(8, 1).Value = 'Gelatin'
$wsSM.Cells.Item(8, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65802'
$wsSM.Cells.Item(9, 1).Value = 'DMSO (serum)'
$wsSM.Cells.Item(9, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125'
$wsSM.Cells.Item(10, 1).Value = 'CMC'
$wsSM.Cells.Item(10, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83594'
$wsSM.Cells.Item(11, 1).Value = '2% PFA/2.5% Glutaraldehyde'
$wsSM.Cells.Item(11, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332'
$wsSM.Cells.Item(12, 1).Value = 'Methanol'
$wsSM.Cells.Item(12, 2).Value = 'http://purl.obolibrary.org/obo/CHEBI_17790'
$wsSM.Cells.Item(13, 1).Value = 'PAXgene tissue kit (PXT)'
$wsSM.Cells.Item(13, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113'
$wsSM.Cells.Item(14, 1).Value = 'PBS'
$wsSM.Cells.Item(14, 2).Value = 'http://purl.obolibrary.org/obo/OBI_0100046'
$wsSM.Cells.Item(15, 1).Value = '1X quench buffer'
$wsSM.Cells.Item(15, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000427'
$wsSM.Cells.Item(16, 1).Value = 'Ethanol'
$wsSM.Cells.Item(16, 2).Value = 'http://purl.obolibrary.org/obo/CHEBI_16236'
$wsSM.Cells.Item(17, 1).Value = 'Formic acid in water'
$wsSM.Cells.Item(17, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83719'
$wsSM.Cells.Item(18, 1).Value = 'HPMC-PVP'
$wsSM.Cells.Item(18, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000386'
$wsSM.Cells.Item(19, 1).Value = 'MACS tissue storage solution'
$wsSM.Cells.Item(19, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105'
$wsSM.Cells.Item(20, 1).Value = 'Tris-EDTA'
$wsSM.Cells.Item(20, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000135'
$wsSM.Cells.Item(21, 1).Value = 'Concentrated quench buffer'
$wsSM.Cells.Item(21, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000391'
$wsSM.Cells.Item(22, 1).Value = 'Cryo-EM'
$wsSM.Cells.Item(22, 2).Value = 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333'
$wsSM.Cells.Item(23, 1).Value = 'RNAlater'
$wsSM.Cells.Item(23, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348'
$wsSM.Cells.Item(24, 1).Value = 'FFPE (Paraffin embedded)'
$wsSM.Cells.Item(24, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C143028'
$wsSM.Cells.Item(25, 1).Value = 'None'
$wsSM.Cells.Item(25, 2).Value = 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132'

# ---- 3. Extend the storage_medium data validation range (K column) ----
$ws.Range("K2:K1001").Validation.Formula1 = '''storage_medium''!$A$1:$A$25'

# ---- 4. Bump pav:createdOn timestamp on the .metadata sheet ----
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Range("C2").Value = "2025-10-16T07:27:11-07:00"
